# Apply crypto price/volume updates per commit "Updated cryptos list" run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.620.02'
$ws.Range('E2').Value = '  -3.14%  '
$ws.Range('D3').Value = '3.180.31'
$ws.Range('E3').Value = '  -4.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.608'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.07%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '3.187.07'
$ws.Range('E9').Value = '  -4.69%  '
$ws.Range('E10').Value = '  -5.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.392'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.27%  '
$ws.Range('D13').Value = '3.724.86'
$ws.Range('E13').Value = '  -5.12%  '
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '64.611.64'
$ws.Range('E15').Value = '  -3.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.57%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000159'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.06%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.185.19'
$ws.Range('E18').Value = '  -4.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '418.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.17'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.205'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.500'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000105'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.996'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.29%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.09'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.39'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E38').Value = '  -5.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.76%  '
$ws.Range('D40').Value = '2.707.96'
$ws.Range('E40').Value = '  -5.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.25'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.32'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.718'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0625'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.97%  '
$ws.Range('E47').Value = '  -3.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '291.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0995'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -14.21%  '
